$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New data: two more login test-scenario rows (3 and 4) plus a "Result"
# column value for every data row. Shared-string table must grow in this
# exact order (Passed, then row 3's scenario name, then row 4's scenario
# name, then the new email, then the new password) to reproduce the target
# uniqueCount/order, so the writes below are sequenced deliberately.
# ---------------------------------------------------------------------------

$ws.Range("D2").Value = "Passed"

$ws.Range("A3").Value = "Test for valid username and invalid password"
$ws.Range("D3").Value = "Passed"

$ws.Range("A4").Value = "Test for invalid username and valid password"
$ws.Range("B4").Value = "haguse@gmail.com"
$ws.Range("D4").Value = "Passed"

$ws.Range("C3").Value = "password124"
$ws.Range("C4").Value = "password123"

# ---------------------------------------------------------------------------
# Hyperlinks: B3:B4 share the existing 234dotus@gmail.com mail-to link
# (matching the display text already used on B2), and B4 additionally gets
# its own personal mail-to link for haguse@gmail.com.
# ---------------------------------------------------------------------------

$ws.Hyperlinks.Add($ws.Range("B3:B4"), "mailto:234dotus@gmail.com", "", "", "234dotus@gmail.com")
$ws.Range("B3:B4").Style = $ws.Range("B2").Style

$ws.Range("B4").Value = "haguse@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:haguse@gmail.com")
$ws.Range("B4").Style = $ws.Range("B2").Style

# Restore the view: scrolled down one row with D8 as the active cell.
$ws.Range("D8").Select()
